$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AMD CPU")

# Row 2
$ws.Range("B2").Value = [double]"8.3999999999999995E-3"
$ws.Range("C2").Value = [double]"6.1999999999999998E-3"
$ws.Range("D2").Value = [double]"6.9999999999999999E-4"
$ws.Range("G2").Value = [double]"0.2147"
$ws.Range("H2").Value = [double]"0.67910000000000004"
$ws.Range("I2").Value = [double]"7.1999999999999995E-2"
$ws.Range("J2").Value = [double]"4.0000000000000002E-4"
$ws.Range("M2").Value = [double]"1.1736E-2"
$ws.Range("N2").Value = [double]"1.4419999999999999E-3"
$ws.Range("O2").Value = [double]"8.2999999999999998E-5"
$ws.Range("R2").Value = [double]"0.52131899999999998"
$ws.Range("S2").Value = [double]"1.0971500000000001"
$ws.Range("T2").Value = [double]"8.0784999999999996E-2"
$ws.Range("U2").Value = [double]"3.3E-4"

# Row 3
$ws.Range("B3").Value = [double]"8.3999999999999995E-3"
$ws.Range("C3").Value = [double]"2.8999999999999998E-3"
$ws.Range("D3").Value = [double]"2.5000000000000001E-3"
$ws.Range("G3").Value = [double]"0.19900000000000001"
$ws.Range("H3").Value = [double]"0.79710000000000003"
$ws.Range("I3").Value = [double]"4.9799999999999997E-2"
$ws.Range("J3").Value = [double]"2.8E-3"
$ws.Range("M3").Value = [double]"4.986E-3"
$ws.Range("N3").Value = [double]"1.763E-3"
$ws.Range("O3").Value = [double]"4.3199999999999998E-4"
$ws.Range("R3").Value = [double]"0.44561600000000001"
$ws.Range("S3").Value = [double]"0.59089499999999995"
$ws.Range("T3").Value = [double]"8.3520999999999998E-2"
$ws.Range("U3").Value = [double]"1.258E-3"

# Row 4
$ws.Range("B4").Value = [double]"8.5900000000000004E-2"
$ws.Range("C4").Value = [double]"1.9E-2"
$ws.Range("D4").Value = [double]"2.8299999999999999E-2"
$ws.Range("G4").Value = [double]"0.19900000000000001"
$ws.Range("H4").Value = [double]"1.0247999999999999"
$ws.Range("I4").Value = [double]"6.7199999999999996E-2"
$ws.Range("J4").Value = [double]"2.5000000000000001E-2"
$ws.Range("M4").Value = [double]"5.6697999999999998E-2"
$ws.Range("N4").Value = [double]"1.7961999999999999E-2"
$ws.Range("O4").Value = [double]"3.875E-3"
$ws.Range("R4").Value = [double]"0.47459200000000001"
$ws.Range("S4").Value = [double]"0.84513199999999999"
$ws.Range("T4").Value = [double]"9.7727999999999995E-2"
$ws.Range("U4").Value = [double]"7.025E-3"

# Row 5
$ws.Range("B5").Value = [double]"0.84370000000000001"
$ws.Range("C5").Value = [double]"0.44240000000000002"
$ws.Range("D5").Value = [double]"0.2477"
$ws.Range("G5").Value = [double]"0.2001"
$ws.Range("H5").Value = [double]"1.2468999999999999"
$ws.Range("I5").Value = [double]"8.5800000000000001E-2"
$ws.Range("J5").Value = [double]"0.24779999999999999"
$ws.Range("M5").Value = [double]"0.54696199999999995"
$ws.Range("N5").Value = [double]"0.1709"
$ws.Range("O5").Value = [double]"3.5851000000000001E-2"
$ws.Range("R5").Value = [double]"0.47682099999999999"
$ws.Range("S5").Value = [double]"0.82023599999999997"
$ws.Range("T5").Value = [double]"0.105322"
$ws.Range("U5").Value = [double]"6.5559000000000006E-2"

# Row 6
$ws.Range("B6").Value = [double]"8.7974999999999994"
$ws.Range("C6").Value = [double]"2.4841000000000002"
$ws.Range("D6").Value = [double]"2.4857999999999998"
$ws.Range("G6").Value = [double]"0.22470000000000001"
$ws.Range("H6").Value = [double]"2.9426999999999999"
$ws.Range("I6").Value = [double]"2.0160999999999998"
$ws.Range("J6").Value = [double]"2.4756999999999998"
$ws.Range("M6").Value = [double]"5.46922"
$ws.Range("N6").Value = [double]"1.97777"
$ws.Range("O6").Value = [double]"0.34337899999999999"
$ws.Range("R6").Value = [double]"0.37850499999999998"
$ws.Range("S6").Value = [double]"2.9758"
$ws.Range("T6").Value = [double]"0.464864"
$ws.Range("U6").Value = [double]"0.69070699999999996"

# Row 7
$ws.Range("B7").Value = [double]"73.321100000000001"
$ws.Range("C7").Value = [double]"21.799900000000001"
$ws.Range("D7").Value = [double]"23.4847"
$ws.Range("G7").Value = [double]"0.47810000000000002"
$ws.Range("H7").Value = [double]"5.0983999999999998"
$ws.Range("I7").Value = [double]"4.6943999999999999"
$ws.Range("J7").Value = [double]"25.093299999999999"
$ws.Range("M7").Value = [double]"46.990099999999998"
$ws.Range("N7").Value = [double]"16.2499"
$ws.Range("O7").Value = [double]"3.2206899999999998"
$ws.Range("R7").Value = [double]"1.00038"
$ws.Range("S7").Value = [double]"6.0666500000000001"
$ws.Range("T7").Value = [double]"5.2098500000000003"
$ws.Range("U7").Value = [double]"6.5240200000000002"

# Row 8
$ws.Range("B8").Value = [double]"624.23699999999997"
$ws.Range("C8").Value = [double]"188.989"
$ws.Range("D8").Value = [double]"209.03200000000001"
$ws.Range("G8").Value = [double]"0.98419999999999996"
$ws.Range("H8").Value = [double]"29.4176"
$ws.Range("I8").Value = [double]"42.903700000000001"
$ws.Range("J8").Value = [double]"225.434"
$ws.Range("M8").Value = [double]"458.85199999999998"
$ws.Range("N8").Value = [double]"167.935"
$ws.Range("O8").Value = [double]"33.293500000000002"
$ws.Range("R8").Value = [double]"2.84823"
$ws.Range("S8").Value = [double]"27.171700000000001"
$ws.Range("T8").Value = [double]"46.458199999999998"
$ws.Range("U8").Value = [double]"65.723500000000001"

# Row 9
$ws.Range("B9").Value = [double]"6345.27"
$ws.Range("C9").Value = [double]"2010.93"
$ws.Range("D9").Value = [double]"2106.0100000000002"
$ws.Range("G9").Value = [double]"7.2793000000000001"
$ws.Range("H9").Value = [double]"260.947"
$ws.Range("I9").Value = [double]"373.25200000000001"
$ws.Range("J9").Value = [double]"2340.81"
$ws.Range("M9").Value = [double]"4632.47"
$ws.Range("N9").Value = [double]"4517.74"
$ws.Range("O9").Value = [double]"328.387"
$ws.Range("R9").Value = [double]"16.939800000000002"
$ws.Range("S9").Value = [double]"235.39599999999999"
$ws.Range("T9").Value = [double]"431.298"
$ws.Range("U9").Value = [double]"647.99900000000002"

$excel.Calculate()

# Update sheet view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("U13").Select()
